# Applies the commit "Baru ditamabahin sheet Variansi dan Standard Deviasi":
#   - Sheet3 is renamed "Simpangan Rata-rata"; its stray header row (C1/D1,
#     shared string "f") is cleared, and a new summary label cell
#     "Simpangan Rata-rata" is written into A17 (merged A17:B17) next to the
#     existing "SR" value in C17.
#   - Two new, empty worksheets are appended after it: "Variansi" and
#     "Simpangan Baku Standard Deviasi" (the latter becomes the active tab).
#   - Selection/active-cell bookkeeping on Sheet2 / Simpangan Rata-rata /
#     Variansi is updated to match.

$wb = $excel.ActiveWorkbook

# ---- Sheet2: just a different remembered selection -------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("M20").Select()

# ---- Sheet3 -> "Simpangan Rata-rata" ----------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

# Clear the leftover "f" / "x" header cells that used to sit above the table.
$ws3.Range("C1:D1").ClearContents()

# Add the new row-17 label, merged across A17:B17, right before the
# existing "SR" value (now in C17).
$ws3.Range("A17").Value = "Simpangan Rata-rata"
$ws3.Range("A17:B17").Merge()
$ws3.Range("A17:B17").HorizontalAlignment = -4108  # xlCenter
$ws3.Range("C17").HorizontalAlignment = -4152      # xlRight

$ws3.Name = "Simpangan Rata-rata"
$ws3.Range("F18").Select()

# ---- New sheet: "Variansi" --------------------------------------------
$wsVariansi = $wb.Worksheets.Add([System.Type]::Missing, $ws3)
$wsVariansi.Name = "Variansi"
$wsVariansi.Range("E10").Select()

# ---- New sheet: "Simpangan Baku Standard Deviasi" ---------------------
$wsSD = $wb.Worksheets.Add([System.Type]::Missing, $wsVariansi)
$wsSD.Name = "Simpangan Baku Standard Deviasi"

# Leave the newly-added "Simpangan Baku Standard Deviasi" sheet active/selected.
$wsSD.Select()
